$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.029.72"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.46%  '

$ws.Range('D3').Value = "'2.015.35"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.72%  '

$ws.Range('D4').Value = "'0.996"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.70%  '

$ws.Range('D5').Value = "'225.69"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.61%  '

$ws.Range('D6').Value = "'0.606"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.55%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('D8').Value = "'54.86"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.22%  '

$ws.Range('E9').Value = '  -2.85%  '

$ws.Range('D10').Value = "'0.0782"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.91%  '

$ws.Range('E11').Value = '  -4.88%  '

$ws.Range('D12').Value = "'2.310.90"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.73%  '

$ws.Range('D13').Value = "'14.12"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.38%  '

$ws.Range('D14').Value = "'20.20"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.70%  '

$ws.Range('D15').Value = "'0.740"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.20%  '

$ws.Range('D16').Value = "'5.11"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.06%  '

$ws.Range('D17').Value = "'2.012.83"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.01%  '

$ws.Range('D18').Value = "'36.980.50"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.36%  '

$ws.Range('D19').Value = "'6.18"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.37%  '

$ws.Range('D20').Value = "'68.70"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.76%  '

$ws.Range('D21').Value = "'0.0₃0814"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.45%  '

$ws.Range('D22').Value = "'223.08"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.74%  '

$ws.Range('E23').Value = '  -0.05%  '

$ws.Range('E24').Value = '  +1.65%  '

$ws.Range('D25').Value = "'2.18"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -6.40%  '

$ws.Range('D26').Value = "'165.90"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.05%  '

$ws.Range('D27').Value = "'9.14"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.77%  '

$ws.Range('E28').Value = '  -0.10%  '

$ws.Range('D29').Value = "'18.70"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.28%  '

$ws.Range('E30').Value = '  -5.12%  '

$ws.Range('E31').Value = '  -4.06%  '

$ws.Range('D32').Value = "'4.53"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.78%  '

$ws.Range('D33').Value = "'0.0611"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.36%  '

$ws.Range('D34').Value = "'4.41"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.04%  '

$ws.Range('D35').Value = "'2.34"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.57%  '

$ws.Range('E36').Value = '  +0.88%  '

$ws.Range('D37').Value = "'0.999"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.11%  '

$ws.Range('D38').Value = "'3.16"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.28%  '

$ws.Range('D39').Value = "'5.27"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.78%  '

$ws.Range('D40').Value = "'1.482.15"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.49%  '

$ws.Range('E41').Value = '  -5.13%  '

$ws.Range('D42').Value = "'94.89"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.42%  '

$ws.Range('E43').Value = '  -4.41%  '

$ws.Range('D44').Value = "'16.29"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.94%  '

$ws.Range('E45').Value = '  -5.06%  '

$ws.Range('D46').Value = "'1.13"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.36%  '

$ws.Range('D47').Value = "'7.21"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.47%  '

$ws.Range('E48').Value = '  -2.67%  '

$ws.Range('E49').Value = '  -1.13%  '

$ws.Range('D50').Value = "'2.197.67"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.77%  '

$ws.Range('B51').Value = 'FTXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D51').Value = "'3.55"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -11.33%  '
